# Update crypto price/volume data per Wed Sep  6 18:00:07 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.818.93"
$ws.Range("E2").Value = "  -0.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.69"
$ws.Range("E3").Value = "  -0.06%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.43%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.84"
$ws.Range("E5").Value = "  -0.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.503"
$ws.Range("E6").Value = "  -0.62%  "

$ws.Range("E7").Value = "  -0.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.259"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0638"
$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.69"
$ws.Range("E10").Value = "  -3.86%  "

$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.655.37"
$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.25"
$ws.Range("E13").Value = "  -0.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.853.84"
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.557"
$ws.Range("E15").Value = "  -1.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0767"
$ws.Range("E16").Value = "  -0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.96"
$ws.Range("E17").Value = "  -0.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.812.98"
$ws.Range("E18").Value = "  -0.31%  "

$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.44"
$ws.Range("E20").Value = "  +1.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.02"
$ws.Range("E21").Value = "  +1.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.94"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.27"
$ws.Range("E23").Value = "  +2.29%  "

$ws.Range("E24").Value = "  -0.31%  "

$ws.Range("E25").Value = "  -1.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.01"
$ws.Range("E26").Value = "  -1.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.121"
$ws.Range("E27").Value = "  -2.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.87"
$ws.Range("E28").Value = "  +0.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.55"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  -0.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0492"
$ws.Range("E31").Value = "  -0.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  +1.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.26"
$ws.Range("E33").Value = "  +0.64%  "

$ws.Range("E34").Value = "  +1.41%  "

$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.902"
$ws.Range("E36").Value = "  -0.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.56"
$ws.Range("E37").Value = "  +0.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.548"
$ws.Range("E38").Value = "  -1.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.113.00"
$ws.Range("E39").Value = "  -1.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0156"
$ws.Range("E40").Value = "  -0.41%  "

$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.57"
$ws.Range("E42").Value = "  +0.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.17"
$ws.Range("E43").Value = "  +1.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.802"
$ws.Range("E44").Value = "  -0.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.763.50"
$ws.Range("E45").Value = "  -0.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0109"
$ws.Range("E46").Value = "  -0.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.21"
$ws.Range("E47").Value = "  -0.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.72"
$ws.Range("E48").Value = "  -0.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.419"
$ws.Range("E49").Value = "  -2.33%  "

$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.37"
$ws.Range("E50").Value = "  +3.01%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0502"
$ws.Range("E51").Value = "  -0.37%  "
